# edit.ps1 — applies the "Add files via upload" revision to the travel-agent
# project-content document:
#   1. Motif bullet: drop the spell-check proofErr wrapper around "Wordpress"
#      (merges three runs into one run of plain text).
#   2. Client Travel Survey bullet: merge the " - " run with the
#      "Name (First/Last), ..." run into a single run.
#   3. Big restructure: move the "Contact" section to after a new "Blog
#      (Entry List)" section (derived from what used to be the start of the
#      "Administrator Log-In" section), re-create the original
#      "Administrator Log-In" section afterwards, move the
#      <w:lastRenderedPageBreak/> from "View New Itinerary Requests" to
#      "Request Information", and move the _GoBack bookmark into the new
#      "Will talk about travel experiences..." paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: Motif bullet — remove proofErr-wrapped "Wordpress" run, merging
# the whole tail of the sentence into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Bootstrap, Font Awesome, Wordpress and others that are available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bootstrap, Font Awesome, Wordpress and others that are available", 2) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: Client Travel Survey bullet — merge " - " and "Name (First/Last),
# ..." into one run, while leaving the preceding "Client Travel Survey" run
# untouched. Done via whole-paragraph XML replacement so the run layout is
# exact.
# ---------------------------------------------------------------------
function Find-ParaIndex($startAt, $text) {
    for ($i = $startAt; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13,[char]7) -eq $text) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndexStartsWith($startAt, $text) {
    for ($i = $startAt; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

$travelSurveyIdx = Find-ParaIndexStartsWith 1 "Client Travel Survey"
$travelSurveyPara = $d.Paragraphs($travelSurveyIdx)

$travelSurveyXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Client Travel Survey</w:t></w:r><w:r><w:t xml:space="preserve"> - Name (First/Last), Date of Birth, Preferred Port of Departure, E-Mail, Phone, Address, Preferences, etc.</w:t></w:r></w:p>
"@

$travelSurveyPara.Range.InsertXML($travelSurveyXml) | Out-Null

# ---------------------------------------------------------------------
# Edit 3: big restructure of the Contact / Blog / Administrator Log-In
# sections. Locate the "Contact" heading paragraph and the first empty
# paragraph following it (the one that carries the _GoBack bookmark);
# replace that whole span with the new section order.
# ---------------------------------------------------------------------
$contactIdx = Find-ParaIndex 1 "Contact"
$emptyIdx = -1
for ($i = $contactIdx; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13,[char]7) -eq "") {
        $emptyIdx = $i
        break
    }
}

$startRange = $d.Paragraphs($contactIdx).Range
$endRange = $d.Paragraphs($emptyIdx).Range
$blockRange = $d.Range($startRange.Start, $endRange.End)

$newSectionXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Blog (Entry List)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>List of Posts (Each Post will have its own separate html path)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Will talk about </w:t></w:r><w:r><w:t>travel experiences, travel tips, etc.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>Search blog posts by keyword</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>List of Months/Years of Blog Posts Available</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Contact</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Name</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Phone</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>E-Mail</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Social Media</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Request Information</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Administrator Log-In</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>Manage Customers (Add, Delete, Update)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>Manage Trip Details</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>Manage Payments</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>View New Itinerary Requests</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:t>Manage Messages (from Clients and System)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr></w:p>
"@

$blockRange.InsertXML($newSectionXml) | Out-Null

Write-Output "Done."
